$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet used to hold a "Tabla1" Excel table (A1:E3) with a No/FullName/
# Email/CurrentAddress/PermanentAddress header row plus two sample rows, and
# a couple of stray formatted-but-empty cells further down (E4, E17).
# Remove the table object entirely (also drops its XML part + autofilter).
$ws.ListObjects.Item(1).Delete()

# Clear every column's content, number formatting and custom width so the
# sheet starts from a blank slate (this also removes the stray E4/E17 cells
# and the old "No/FullName/..." header row + sample data rows).
$ws.Columns("A:E").Delete()

# New sheet content: a single header-less row holding a login e-mail and its
# password, used by the Selenium data-driven test.
$ws.Range("A1").Value = "intern6@agilethought.com"
$ws.Range("B1").Value = "P@ssw0rd"

# Column A is widened to comfortably fit the e-mail address.
$ws.Columns.Item(1).ColumnWidth = 26
